$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.236.18'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '1.562.68'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.20'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.91%  '
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0870'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '1.785.82'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '1.570.10'
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.518'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").Value = '27.209.55'
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").Value = '0.0₃0701'
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("E22").Value = '  +0.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +1.84%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("E30").Value = '  +1.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0471'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("D33").Value = '1.450.81'
$ws.Range("E33").Value = '  +2.36%  '
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("E35").Value = '  +5.18%  '
$ws.Range("E36").Value = '  +1.01%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("E38").Value = '  +0.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.541'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.812'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.989'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").Value = '1.699.28'
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("E49").Value = '  +2.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0525'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0948'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.95%  '
